# "massive changes with method formatting"
#
# The document has 14 paragraphs:
#   1  Ingredients\n
#   2  <tab>Bell Peppers
#   3  \n
#   4  Instructions\n
#   5  <tab>Wash Bell Peppers thoroughly.\n
#   6  <tab>\n
#   7  (ind 720) Slice Bell Peppers into medium ½ inch pieces.
#   8  (ind 720) \n
#   9  (ind 720) No need to blanch.\n
#   10 (ind 720) \n
#   11 (ind 720) Space out and place on oven pan with parchment paper and dry for approximately 9 hours at 125 degrees.
#   12 (ind 720) \n
#   13 (ind 720) Pour dehydrated pepper pieces into jars. \n
#   14 (ind 720) \n
#
# Target layout:
#   1  Ingredients\n
#   2  16-spaces + Bell Peppers + \n   (tab replaced by 16 spaces; merged with the old
#                                       standalone "\n" paragraph that used to follow it)
#   3  \n                              (new standalone blank-line paragraph, re-created
#                                       right before "Instructions")
#   4  Instructions\n
#   5  Wash Bell Peppers thoroughly.\n (leading tab removed)
#   6  <tab>\n                        (unchanged)
#   7  Slice Bell Peppers into medium ½ inch pieces.   (left indent removed)
#   8  (ind 720) \n                                    (unchanged)
#   9  No need to blanch.\n                            (left indent removed)
#   10 (ind 720) \n                                    (unchanged)
#   11 Space out and place on oven pan ... degrees.     (left indent removed)
#   12 (ind 720) \n                                    (unchanged)
#   13 Pour dehydrated pepper pieces into jars. \n      (left indent removed)
#   14 (ind 720) \n                                    (unchanged)

$d = $word.ActiveDocument

# --- 1. Merge the "Bell Peppers" paragraph with the blank "\n" paragraph after it ---
# (deleting the paragraph mark at the end of paragraph 2 pulls paragraph 3's run into it)
$pBellPeppers = $d.Paragraphs.Item(2)
$rBellPeppers = $pBellPeppers.Range
$null = $d.Range($rBellPeppers.End - 1, $rBellPeppers.End).Delete()

# Replace the leading tab in that paragraph with 16 spaces.
$pBellPeppers = $d.Paragraphs.Item(2)
$null = $pBellPeppers.Range.Find.Execute("^t", $false, $false, $false, $false, $false, $true, 1, $false, "                ", 2)

# --- 2. Insert a fresh blank "\n" paragraph right before "Instructions" ---
$pInstructions = $d.Paragraphs.Item(3)
$null = $pInstructions.Range.InsertParagraphBefore()
$pNewBlank = $d.Paragraphs.Item(3)
$pNewBlank.Range.Text = "\n"

# --- 3. Remove the leading tab on "Wash Bell Peppers thoroughly." ---
$pWash = $d.Paragraphs.Item(5)
$null = $pWash.Range.Find.Execute("^t", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- 4. Drop the 720-twip left indent from the four instructional text paragraphs ---
$pSlice = $d.Paragraphs.Item(7)
$pSlice.Range.ParagraphFormat.LeftIndent = 0

$pBlanch = $d.Paragraphs.Item(9)
$pBlanch.Range.ParagraphFormat.LeftIndent = 0

$pSpaceOut = $d.Paragraphs.Item(11)
$pSpaceOut.Range.ParagraphFormat.LeftIndent = 0

$pPour = $d.Paragraphs.Item(13)
$pPour.Range.ParagraphFormat.LeftIndent = 0

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
